# Apply the price/volume refresh captured in the commit's OOXML diff.
# Column D ("Price") values are forced to text (leading apostrophe, like typing
# '123.45 into Excel) so numeric-looking strings keep their exact original
# formatting (e.g. trailing zeros such as "0.150" or multi-dot figures such as
# "53.938.42") instead of being silently parsed into floating point numbers.
# Column E ("Volume(1h)") values are plain padded percentage strings and are
# never mistaken for numbers, so they are written as-is.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'" + "54.263.71"
$ws.Range("E2").Value = "  -3.64%  "
$ws.Range("D3").Value = "'" + "2.262.70"
$ws.Range("E3").Value = "  -4.58%  "
$ws.Range("D4").Value = "'" + "0.998"
$ws.Range("E4").Value = "  -0.28%  "
$ws.Range("D5").Value = "'" + "491.25"
$ws.Range("D6").Value = "'" + "126.94"
$ws.Range("E6").Value = "  -2.39%  "
$ws.Range("D7").Value = "'" + "0.998"
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").Value = "  -3.35%  "
$ws.Range("D9").Value = "'" + "2.262.44"
$ws.Range("E9").Value = "  -4.96%  "
$ws.Range("D10").Value = "'" + "0.0934"
$ws.Range("E10").Value = "  -5.28%  "
$ws.Range("D11").Value = "'" + "0.150"
$ws.Range("E11").Value = "  -0.02%  "
$ws.Range("E12").Value = "  -1.86%  "
$ws.Range("E13").Value = "  -5.00%  "
$ws.Range("D14").Value = "'" + "2.659.65"
$ws.Range("E14").Value = "  -4.76%  "
$ws.Range("D15").Value = "'" + "21.46"
$ws.Range("E15").Value = "  -0.21%  "
$ws.Range("D16").Value = "'" + "53.938.42"
$ws.Range("E16").Value = "  -4.17%  "
$ws.Range("E17").Value = "  -3.00%  "
$ws.Range("D18").Value = "'" + "2.245.13"
$ws.Range("E18").Value = "  -5.08%  "
$ws.Range("E19").Value = "  -2.64%  "
$ws.Range("E20").Value = "  -0.08%  "
$ws.Range("D21").Value = "'" + "297.80"
$ws.Range("E21").Value = "  -3.71%  "
$ws.Range("D22").Value = "'" + "6.26"
$ws.Range("E22").Value = "  -0.57%  "
$ws.Range("D23").Value = "'" + "0.999"
$ws.Range("E23").Value = "  -0.02%  "
$ws.Range("D24").Value = "'" + "63.70"
$ws.Range("E24").Value = "  -3.38%  "
$ws.Range("D25").Value = "'" + "0.999"
$ws.Range("E25").Value = "  +0.20%  "
$ws.Range("E26").Value = "  +0.54%  "
$ws.Range("E27").Value = "  -0.96%  "
$ws.Range("D28").Value = "'" + "2.327.22"
$ws.Range("E28").Value = "  -6.25%  "
$ws.Range("E29").Value = "  -1.50%  "
$ws.Range("D30").Value = "'" + "163.22"
$ws.Range("E30").Value = "  -5.69%  "
$ws.Range("E31").Value = "  -3.06%  "
$ws.Range("D32").Value = "'" + "0.0₃0680"
$ws.Range("E32").Value = "  -4.43%  "
$ws.Range("E34").Value = "  -0.56%  "
$ws.Range("D35").Value = "'" + "0.999"
$ws.Range("E35").Value = "  +0.25%  "
$ws.Range("E36").Value = "  +0.15%  "
$ws.Range("D37").Value = "'" + "17.41"
$ws.Range("E37").Value = "  -1.13%  "
$ws.Range("E38").Value = "  +0.56%  "
$ws.Range("E39").Value = "  +1.47%  "
$ws.Range("E40").Value = "  -1.92%  "
$ws.Range("D41").Value = "'" + "35.35"
$ws.Range("E41").Value = "  -2.89%  "
$ws.Range("E42").Value = "  +0.64%  "
$ws.Range("E43").Value = "  +0.47%  "
$ws.Range("E44").Value = "  -1.46%  "
$ws.Range("D45").Value = "'" + "126.94"
$ws.Range("E45").Value = "  +0.90%  "
$ws.Range("D46").Value = "'" + "4.79"
$ws.Range("E46").Value = "  -2.85%  "
$ws.Range("E47").Value = "  -0.84%  "
$ws.Range("D48").Value = "'" + "242.24"
$ws.Range("E48").Value = "  +1.75%  "
$ws.Range("E49").Value = "  -3.63%  "
$ws.Range("E50").Value = "  -0.65%  "
$ws.Range("E51").Value = "  -1.90%  "
